# Apply "Demais comandos linux, diretório e permissões." edit.
$wb = $excel.ActiveWorkbook

$wsComandos = $wb.Worksheets.Item("Comandos")
$wsGrupos = $wb.Worksheets.Item("GrupoDeComanodos")

# ---------------------------------------------------------------------------
# Sheet "Comandos" - new rows 33..49
# ---------------------------------------------------------------------------

# Row 33
$wsComandos.Range("A33").Value = "apt"
$wsComandos.Range("B33").Value = "Instalação"

# Row 34
$wsComandos.Range("A34").Value = "sudo"

# Row 35
$wsComandos.Range("A35").Value = "apt-get"

# Row 36
$wsComandos.Range("A36").Value = "ls -l"

# Row 37
$wsComandos.Range("A37").Value = "chmod [opções] modo arquivo"

# Row 38
$wsComandos.Range("A38").Value = "addgroup"
$wsComandos.Range("B38").Value = "Grupos de Usuários"

# Row 39
$wsComandos.Range("A39").Value = "adduser [usuario] [grupo]"
$wsComandos.Range("B39").Value = "Ususários"

# Row 40 (taller row with wrapped description text)
$wsComandos.Range("A40").Value = "id"
$wsComandos.Range("B40").Value = "Ususários"
$wsComandos.Range("C40").Value = "Utilize o comando id sempre que quiser obter informações sobre um usuário`ndo sistema. Ele irá retornar UID, GID e os grupos aos quais o usuário pertence."
$wsComandos.Range("C40").WrapText = $true
$wsComandos.Rows.Item(40).RowHeight = 60

# Row 41
$wsComandos.Range("A41").Value = "groups"
$wsComandos.Range("B41").Value = "Ususários"

# Row 42
$wsComandos.Range("A42").Value = " delUser [usuario] [grupo]"

# Row 43
$wsComandos.Range("A43").Value = "chown"

# Row 44
$wsComandos.Range("A44").Value = "chgrp [grupo] [arquivo]"

# Row 45
$wsComandos.Range("A45").Value = "aptitude"

# Row 46
$wsComandos.Range("A46").Value = "add-apt-repository"

# Row 47
$wsComandos.Range("A47").Value = "sudo apt-get install build-essential"

# Row 48
$wsComandos.Range("A48").Value = "sudo apt-get install apache2"

# Row 49
$wsComandos.Range("A49").Value = "apt-cache"
$wsComandos.Range("C49").Value = "Sempre que quiser obter informações sobre um pacote"

# ---------------------------------------------------------------------------
# Sheet "GrupoDeComanodos" - new rows 9..10
# ---------------------------------------------------------------------------
$wsGrupos.Range("A9").Value = "Instalação"
$wsGrupos.Range("A10").Value = "Grupos de Usuários"

# ---------------------------------------------------------------------------
# Defined name "Grupos" now covers the two extra rows
# ---------------------------------------------------------------------------
$wb.Names.Item("Grupos").RefersTo = "=GrupoDeComanodos!`$A`$1:`$A`$10"

# ---------------------------------------------------------------------------
# Update selections to mirror the author's final cursor position
# ---------------------------------------------------------------------------
[void]$wsGrupos.Range("A11").Select()
[void]$wsComandos.Activate()
[void]$wsComandos.Range("C50").Select()

Write-Output "done"
